# Chapter 5 Exercises edit script
$d = $word.ActiveDocument

# 1) Remove the "_GoBack" bookmark from the "True. In a partially ..." paragraph
#    (it will be re-added around the picture paragraph below). Doing this first
#    frees up bookmark id 0 so the new bookmark reuses the same id.
$pTrue = $d.Paragraphs(9)
$pTrue.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>True.  In a partially observable game, it does not provide much information to know what move the second player will make.</w:t></w:r></w:p>')

# 2) Title: "Homework Chapter 5" -> "Chapter " + "5" + " Exercises" (3 runs)
$pTitle = $d.Paragraphs(1)
$pTitle.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Title"/></w:pPr><w:r><w:t xml:space="preserve">Chapter </w:t></w:r><w:r><w:t>5</w:t></w:r><w:r><w:t xml:space="preserve"> Exercises</w:t></w:r></w:p>')

# 3) Drop the lastRenderedPageBreak on the "Exercise 5.8" heading
$pEx58 = $d.Paragraphs(2)
$pEx58.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Exercise</w:t></w:r><w:r w:rsidR="009947D2"><w:t xml:space="preserve"> 5.8</w:t></w:r></w:p>')

# 4) Picture paragraph: resize image, stamp anchorId/editId, add _GoBack bookmark
#    after the picture run, and append a new blank paragraph right after it.
$pPic = $d.Paragraphs(3)
$pPic.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:rPr><w:noProof/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="490E1AD9" wp14:editId="3D28D757"><wp:extent cx="5486400" cy="6108700"/><wp:effectExtent l="0" t="0" r="0" b="12700"/><wp:docPr id="1" name="Picture 1"/><wp:cNvGraphicFramePr><a:graphicFrameLocks noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr id="0" name="thumbnail_image1.jpg"/><pic:cNvPicPr/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId5"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="5486400" cy="6108700"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>')

# 5) "Exercise 5.21" heading now picks up the lastRenderedPageBreak that used to
#    sit on "Exercise 5.8" (paragraph indices shifted by the new blank paragraph).
$pEx521 = $d.Paragraphs(5)
$pEx521.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Exercise 5.21</w:t></w:r></w:p>')

Write-Host "Paragraphs:" $d.Paragraphs.Count
